$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from row 6 (same category/country/product as the new row)
# down to row 11 so the new row's date cell (A11) picks up the same
# date number-format style (s="2") instead of minting a new style.
$ws.Range("A6:M6").Copy($ws.Range("A11:M11"))

# Populate the new row's values (mirrors existing row 6 except the date).
$ws.Range("A11").Value = 45026
$ws.Range("B11").Value = 394
$ws.Range("C11").Value = 5
$ws.Range("D11").Value = "Automitivos"
$ws.Range("E11").Value = 639
$ws.Range("F11").Value = "USA"
$ws.Range("G11").Value = 8
$ws.Range("H11").Value = "Limpa vidros"
$ws.Range("I11").Value = 250
$ws.Range("J11").Value = 290
$ws.Range("K11").Value = 300
$ws.Range("L11").Value = 87000
$ws.Range("M11").Value = 19392

# Match the selection shown in the saved file.
$ws.Range("L11").Select()
